$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 42
$ws.Range("B42").Value = 6078859
$ws.Range("F42").Value = "Union La Calera"
$ws.Range("G42").Value = "Nublense"
$ws.Range("H42").Value = 1
$ws.Range("I42").Value = 1
$ws.Range("J42").Value = "D"
$ws.Range("K42").Value = 2.55
$ws.Range("L42").Value = 3.4
$ws.Range("M42").Value = 2.7
$ws.Range("N42").Value = 1.95
$ws.Range("O42").Value = 3.6
$ws.Range("P42").Value = 3.8
$ws.Range("Q42").Value = -0.5
$ws.Range("R42").Value = 1.95
$ws.Range("S42").Value = 1.85
$ws.Range("T42").Value = 2.5
$ws.Range("U42").Value = 1.925
$ws.Range("V42").Value = 1.875
$ws.Range("W42").Value = -1
$ws.Range("X42").Value = 2.6
$ws.Range("Y42").Value = -1
$ws.Range("Z42").Value = -1
$ws.Range("AA42").Value = 0.8500000000000001
$ws.Range("AB42").Value = -1
$ws.Range("AC42").Value = 0.875

# Row 43
$ws.Range("B43").Value = 6074216
$ws.Range("F43").Value = "Magallanes"
$ws.Range("G43").Value = "Cobresal"
$ws.Range("H43").Value = 2
$ws.Range("I43").Value = 1
$ws.Range("J43").Value = "H"
$ws.Range("K43").Value = 2.8
$ws.Range("L43").Value = 3.4
$ws.Range("M43").Value = 2.45
$ws.Range("N43").Value = 2.55
$ws.Range("O43").Value = 3.4
$ws.Range("P43").Value = 2.625
$ws.Range("Q43").Value = 0
$ws.Range("R43").Value = 1.875
$ws.Range("S43").Value = 1.925
$ws.Range("T43").Value = 2.5
$ws.Range("U43").Value = 1.85
$ws.Range("V43").Value = 1.95
$ws.Range("W43").Value = 1.55
$ws.Range("X43").Value = -1
$ws.Range("Y43").Value = -1
$ws.Range("Z43").Value = 0.875
$ws.Range("AA43").Value = -1
$ws.Range("AB43").Value = 0.8500000000000001
$ws.Range("AC43").Value = -1

# Row 45
$ws.Range("B45").Value = 6078935
$ws.Range("F45").Value = "Everton de Vina"
$ws.Range("G45").Value = "Audax Italiano"
$ws.Range("H45").Value = 3
$ws.Range("I45").Value = 1
$ws.Range("J45").Value = "H"
$ws.Range("K45").Value = 1.95
$ws.Range("L45").Value = 3.5
$ws.Range("M45").Value = 3.8
$ws.Range("N45").Value = 1.75
$ws.Range("O45").Value = 4
$ws.Range("P45").Value = 4.5
$ws.Range("Q45").Value = -0.75
$ws.Range("R45").Value = 1.95
$ws.Range("S45").Value = 1.85
$ws.Range("T45").Value = 2.5
$ws.Range("U45").Value = 1.925
$ws.Range("V45").Value = 1.875
$ws.Range("W45").Value = 0.75
$ws.Range("X45").Value = -1
$ws.Range("Y45").Value = -1
$ws.Range("Z45").Value = 0.95
$ws.Range("AA45").Value = -1
$ws.Range("AB45").Value = 0.925
$ws.Range("AC45").Value = -1

# Row 46
$ws.Range("B46").Value = 6073148
$ws.Range("F46").Value = "Deportes Copiapo"
$ws.Range("G46").Value = "OHiggins"
$ws.Range("H46").Value = 2
$ws.Range("I46").Value = 2
$ws.Range("J46").Value = "D"
$ws.Range("K46").Value = 3.2
$ws.Range("L46").Value = 3.4
$ws.Range("M46").Value = 2.2
$ws.Range("N46").Value = 2.4
$ws.Range("O46").Value = 3.25
$ws.Range("P46").Value = 3
$ws.Range("Q46").Value = -0.25
$ws.Range("R46").Value = 2.05
$ws.Range("S46").Value = 1.75
$ws.Range("T46").Value = 2.25
$ws.Range("U46").Value = 1.8
$ws.Range("V46").Value = 2
$ws.Range("W46").Value = -1
$ws.Range("X46").Value = 2.25
$ws.Range("Y46").Value = -1
$ws.Range("Z46").Value = -0.5
$ws.Range("AA46").Value = 0.375
$ws.Range("AB46").Value = 0.8
$ws.Range("AC46").Value = -1

# Row 66
$ws.Range("B66").Value = 6075784
$ws.Range("F66").Value = "Magallanes"
$ws.Range("G66").Value = "Universidad Catolica"
$ws.Range("H66").Value = 1
$ws.Range("I66").Value = 1
$ws.Range("J66").Value = "D"
$ws.Range("K66").Value = 3
$ws.Range("L66").Value = 3.3
$ws.Range("M66").Value = 2.15
$ws.Range("N66").Value = 3.4
$ws.Range("O66").Value = 3.6
$ws.Range("P66").Value = 2.05
$ws.Range("Q66").Value = 0.25
$ws.Range("R66").Value = 2.025
$ws.Range("S66").Value = 1.825
$ws.Range("T66").Value = 2.75
$ws.Range("U66").Value = 1.925
$ws.Range("V66").Value = 1.925
$ws.Range("W66").Value = -1
$ws.Range("X66").Value = 2.6
$ws.Range("Y66").Value = -1
$ws.Range("Z66").Value = 0.5125
$ws.Range("AA66").Value = -0.5
$ws.Range("AB66").Value = -1
$ws.Range("AC66").Value = 0.925

# Row 67
$ws.Range("B67").Value = 6078866
$ws.Range("F67").Value = "Union Espanola"
$ws.Range("G67").Value = "Palestino"
$ws.Range("H67").Value = 2
$ws.Range("I67").Value = 0
$ws.Range("J67").Value = "H"
$ws.Range("K67").Value = 2.625
$ws.Range("L67").Value = 3.3
$ws.Range("M67").Value = 2.375
$ws.Range("N67").Value = 3.3
$ws.Range("O67").Value = 3.5
$ws.Range("P67").Value = 2.15
$ws.Range("Q67").Value = 0.25
$ws.Range("R67").Value = 1.95
$ws.Range("S67").Value = 1.85
$ws.Range("T67").Value = 2.5
$ws.Range("U67").Value = 1.875
$ws.Range("V67").Value = 1.925
$ws.Range("W67").Value = 2.3
$ws.Range("X67").Value = -1
$ws.Range("Y67").Value = -1
$ws.Range("Z67").Value = 0.95
$ws.Range("AA67").Value = -1
$ws.Range("AB67").Value = -1
$ws.Range("AC67").Value = 0.925

# Row 82
$ws.Range("B82").Value = 6078871
$ws.Range("F82").Value = "Cobresal"
$ws.Range("G82").Value = "Curico Unido"
$ws.Range("H82").Value = 2
$ws.Range("I82").Value = 0
$ws.Range("J82").Value = "H"
$ws.Range("K82").Value = 1.75
$ws.Range("L82").Value = 3.6
$ws.Range("M82").Value = 4.2
$ws.Range("N82").Value = 1.444
$ws.Range("O82").Value = 4.75
$ws.Range("P82").Value = 7
$ws.Range("Q82").Value = -1.25
$ws.Range("R82").Value = 1.95
$ws.Range("S82").Value = 1.85
$ws.Range("T82").Value = 3
$ws.Range("U82").Value = 1.825
$ws.Range("V82").Value = 1.975
$ws.Range("W82").Value = 0.444
$ws.Range("X82").Value = -1
$ws.Range("Y82").Value = -1
$ws.Range("Z82").Value = 0.95
$ws.Range("AA82").Value = -1
$ws.Range("AB82").Value = -1
$ws.Range("AC82").Value = 0.9750000000000001

# Row 83
$ws.Range("B83").Value = 6075788
$ws.Range("F83").Value = "Magallanes"
$ws.Range("G83").Value = "Deportes Copiapo"
$ws.Range("H83").Value = 0
$ws.Range("I83").Value = 2
$ws.Range("J83").Value = "A"
$ws.Range("K83").Value = 1.95
$ws.Range("L83").Value = 3.4
$ws.Range("M83").Value = 3.5
$ws.Range("N83").Value = 2
$ws.Range("O83").Value = 3.5
$ws.Range("P83").Value = 3.75
$ws.Range("Q83").Value = -0.5
$ws.Range("R83").Value = 2.025
$ws.Range("S83").Value = 1.825
$ws.Range("T83").Value = 2.5
$ws.Range("U83").Value = 1.875
$ws.Range("V83").Value = 1.975
$ws.Range("W83").Value = -1
$ws.Range("X83").Value = -1
$ws.Range("Y83").Value = 2.75
$ws.Range("Z83").Value = -1
$ws.Range("AA83").Value = 0.825
$ws.Range("AB83").Value = -1
$ws.Range("AC83").Value = 0.9750000000000001

# Row 217
$ws.Range("B217").Value = 7494647
$ws.Range("F217").Value = "Huachipato"
$ws.Range("G217").Value = "Universidad Catolica"
$ws.Range("H217").Value = 1
$ws.Range("I217").Value = 1
$ws.Range("J217").Value = "D"
$ws.Range("K217").Value = 2.2
$ws.Range("L217").Value = 3.4
$ws.Range("M217").Value = 3.2
$ws.Range("N217").Value = 1.8
$ws.Range("O217").Value = 3.6
$ws.Range("P217").Value = 4.333
$ws.Range("Q217").Value = -0.75
$ws.Range("R217").Value = 1.975
$ws.Range("S217").Value = 1.875
$ws.Range("T217").Value = 2.75
$ws.Range("U217").Value = 1.975
$ws.Range("V217").Value = 1.875
$ws.Range("W217").Value = -1
$ws.Range("X217").Value = 2.6
$ws.Range("Y217").Value = -1
$ws.Range("Z217").Value = -1
$ws.Range("AA217").Value = 0.875
$ws.Range("AB217").Value = -1
$ws.Range("AC217").Value = 0.875

# Row 218
$ws.Range("B218").Value = 7494646
$ws.Range("F218").Value = "OHiggins"
$ws.Range("G218").Value = "Cobresal"
$ws.Range("H218").Value = 0
$ws.Range("I218").Value = 0
$ws.Range("J218").Value = "D"
$ws.Range("K218").Value = 3
$ws.Range("L218").Value = 3.4
$ws.Range("M218").Value = 2.3
$ws.Range("N218").Value = 2.1
$ws.Range("O218").Value = 3.5
$ws.Range("P218").Value = 3.5
$ws.Range("Q218").Value = -0.25
$ws.Range("R218").Value = 1.8
$ws.Range("S218").Value = 2.05
$ws.Range("T218").Value = 2.75
$ws.Range("U218").Value = 1.975
$ws.Range("V218").Value = 1.875
$ws.Range("W218").Value = -1
$ws.Range("X218").Value = 2.5
$ws.Range("Y218").Value = -1
$ws.Range("Z218").Value = -0.5
$ws.Range("AA218").Value = 0.5249999999999999
$ws.Range("AB218").Value = -1
$ws.Range("AC218").Value = 0.875

# Row 220
$ws.Range("B220").Value = 6077497
$ws.Range("F220").Value = "Deportes Copiapo"
$ws.Range("G220").Value = "Nublense"
$ws.Range("H220").Value = 1
$ws.Range("I220").Value = 1
$ws.Range("J220").Value = "D"
$ws.Range("K220").Value = 2.6
$ws.Range("L220").Value = 3.4
$ws.Range("M220").Value = 2.6
$ws.Range("N220").Value = 2.8
$ws.Range("O220").Value = 3.2
$ws.Range("P220").Value = 2.7
$ws.Range("Q220").Value = 0
$ws.Range("R220").Value = 1.95
$ws.Range("S220").Value = 1.9
$ws.Range("T220").Value = 2.25
$ws.Range("U220").Value = 2
$ws.Range("V220").Value = 1.85
$ws.Range("W220").Value = -1
$ws.Range("X220").Value = 2.2
$ws.Range("Y220").Value = -1
$ws.Range("Z220").Value = 0
$ws.Range("AA220").Value = -0
$ws.Range("AB220").Value = -0.5
$ws.Range("AC220").Value = 0.425

# Row 221
$ws.Range("B221").Value = 6077763
$ws.Range("F221").Value = "Curico Unido"
$ws.Range("G221").Value = "Magallanes"
$ws.Range("H221").Value = 3
$ws.Range("I221").Value = 4
$ws.Range("J221").Value = "A"
$ws.Range("K221").Value = 2.15
$ws.Range("L221").Value = 3.5
$ws.Range("M221").Value = 3.2
$ws.Range("N221").Value = 2.625
$ws.Range("O221").Value = 3.5
$ws.Range("P221").Value = 2.6
$ws.Range("Q221").Value = 0
$ws.Range("R221").Value = 1.95
$ws.Range("S221").Value = 1.9
$ws.Range("T221").Value = 2.75
$ws.Range("U221").Value = 1.975
$ws.Range("V221").Value = 1.875
$ws.Range("W221").Value = -1
$ws.Range("X221").Value = -1
$ws.Range("Y221").Value = 1.6
$ws.Range("Z221").Value = -1
$ws.Range("AA221").Value = 0.8999999999999999
$ws.Range("AB221").Value = 0.9750000000000001
$ws.Range("AC221").Value = -1

# Row 223
$ws.Range("B223").Value = 6077498
$ws.Range("F223").Value = "Universidad Catolica"
$ws.Range("G223").Value = "Deportes Copiapo"
$ws.Range("H223").Value = 2
$ws.Range("I223").Value = 2
$ws.Range("J223").Value = "D"
$ws.Range("K223").Value = 1.65
$ws.Range("L223").Value = 3.8
$ws.Range("M223").Value = 5.25
$ws.Range("N223").Value = 1.909
$ws.Range("O223").Value = 3.6
$ws.Range("P223").Value = 4.2
$ws.Range("Q223").Value = -0.5
$ws.Range("R223").Value = 1.85
$ws.Range("S223").Value = 2
$ws.Range("T223").Value = 2.75
$ws.Range("U223").Value = 2.025
$ws.Range("V223").Value = 1.825
$ws.Range("W223").Value = -1
$ws.Range("X223").Value = 2.6
$ws.Range("Y223").Value = -1
$ws.Range("Z223").Value = -1
$ws.Range("AA223").Value = 1
$ws.Range("AB223").Value = 1.025
$ws.Range("AC223").Value = -1

# Row 224
$ws.Range("B224").Value = 6078266
$ws.Range("F224").Value = "Palestino"
$ws.Range("G224").Value = "Curico Unido"
$ws.Range("H224").Value = 4
$ws.Range("I224").Value = 0
$ws.Range("J224").Value = "H"
$ws.Range("K224").Value = 1.533
$ws.Range("L224").Value = 4
$ws.Range("M224").Value = 6
$ws.Range("N224").Value = 1.363
$ws.Range("O224").Value = 4.75
$ws.Range("P224").Value = 7.5
$ws.Range("Q224").Value = -1.5
$ws.Range("R224").Value = 2.025
$ws.Range("S224").Value = 1.825
$ws.Range("T224").Value = 3
$ws.Range("U224").Value = 1.9
$ws.Range("V224").Value = 1.95
$ws.Range("W224").Value = 0.363
$ws.Range("X224").Value = -1
$ws.Range("Y224").Value = -1
$ws.Range("Z224").Value = 1.025
$ws.Range("AA224").Value = -1
$ws.Range("AB224").Value = 0.8999999999999999
$ws.Range("AC224").Value = -1

# Row 225
$ws.Range("B225").Value = 6078265
$ws.Range("F225").Value = "Audax Italiano"
$ws.Range("G225").Value = "Magallanes"
$ws.Range("H225").Value = 0
$ws.Range("I225").Value = 2
$ws.Range("J225").Value = "A"
$ws.Range("K225").Value = 1.666
$ws.Range("L225").Value = 3.75
$ws.Range("M225").Value = 5
$ws.Range("N225").Value = 2.25
$ws.Range("O225").Value = 3.3
$ws.Range("P225").Value = 3.3
$ws.Range("Q225").Value = -0.25
$ws.Range("R225").Value = 1.95
$ws.Range("S225").Value = 1.85
$ws.Range("T225").Value = 2.5
$ws.Range("U225").Value = 1.8
$ws.Range("V225").Value = 2
$ws.Range("W225").Value = -1
$ws.Range("X225").Value = -1
$ws.Range("Y225").Value = 2.3
$ws.Range("Z225").Value = -1
$ws.Range("AA225").Value = 0.8500000000000001
$ws.Range("AB225").Value = -1
$ws.Range("AC225").Value = 1

# Row 230
$ws.Range("B230").Value = 6078267
$ws.Range("F230").Value = "Huachipato"
$ws.Range("G230").Value = "Audax Italiano"
$ws.Range("H230").Value = 2
$ws.Range("I230").Value = 0
$ws.Range("J230").Value = "H"
$ws.Range("K230").Value = 1.5
$ws.Range("L230").Value = 4.333
$ws.Range("M230").Value = 6
$ws.Range("N230").Value = 1.444
$ws.Range("O230").Value = 4.75
$ws.Range("P230").Value = 7
$ws.Range("Q230").Value = -1.25
$ws.Range("R230").Value = 2.025
$ws.Range("S230").Value = 1.825
$ws.Range("T230").Value = 2.75
$ws.Range("U230").Value = 1.8
$ws.Range("V230").Value = 2.05
$ws.Range("W230").Value = 0.444
$ws.Range("X230").Value = -1
$ws.Range("Y230").Value = -1
$ws.Range("Z230").Value = 1.025
$ws.Range("AA230").Value = -1
$ws.Range("AB230").Value = -1
$ws.Range("AC230").Value = 1.05

# Row 231
$ws.Range("B231").Value = 6143704
$ws.Range("F231").Value = "Curico Unido"
$ws.Range("G231").Value = "Colo Colo"
$ws.Range("H231").Value = 0
$ws.Range("I231").Value = 1
$ws.Range("J231").Value = "A"
$ws.Range("K231").Value = 6.5
$ws.Range("L231").Value = 4.75
$ws.Range("M231").Value = 1.4
$ws.Range("N231").Value = 12
$ws.Range("O231").Value = 8.5
$ws.Range("P231").Value = 1.166
$ws.Range("Q231").Value = 2
$ws.Range("R231").Value = 2
$ws.Range("S231").Value = 1.8
$ws.Range("T231").Value = 3.25
$ws.Range("U231").Value = 1.875
$ws.Range("V231").Value = 1.925
$ws.Range("W231").Value = -1
$ws.Range("X231").Value = -1
$ws.Range("Y231").Value = 0.1659999999999999
$ws.Range("Z231").Value = 1
$ws.Range("AA231").Value = -1
$ws.Range("AB231").Value = -1
$ws.Range("AC231").Value = 0.925

# Row 232
$ws.Range("B232").Value = 6078997
$ws.Range("F232").Value = "Union Espanola"
$ws.Range("G232").Value = "Cobresal"
$ws.Range("H232").Value = 1
$ws.Range("I232").Value = 0
$ws.Range("J232").Value = "H"
$ws.Range("K232").Value = 3.8
$ws.Range("L232").Value = 3.6
$ws.Range("M232").Value = 1.909
$ws.Range("N232").Value = 2.7
$ws.Range("O232").Value = 3.6
$ws.Range("P232").Value = 2.45
$ws.Range("Q232").Value = 0
$ws.Range("R232").Value = 1.975
$ws.Range("S232").Value = 1.825
$ws.Range("T232").Value = 2.75
$ws.Range("U232").Value = 1.775
$ws.Range("V232").Value = 2.025
$ws.Range("W232").Value = 1.7
$ws.Range("X232").Value = -1
$ws.Range("Y232").Value = -1
$ws.Range("Z232").Value = 0.9750000000000001
$ws.Range("AA232").Value = -1
$ws.Range("AB232").Value = -1
$ws.Range("AC232").Value = 1.025

# Row 233
$ws.Range("B233").Value = 6077499
$ws.Range("F233").Value = "Deportes Copiapo"
$ws.Range("G233").Value = "Everton de Vina"
$ws.Range("H233").Value = 2
$ws.Range("I233").Value = 0
$ws.Range("J233").Value = "H"
$ws.Range("K233").Value = 2.1
$ws.Range("L233").Value = 3.4
$ws.Range("M233").Value = 3.4
$ws.Range("N233").Value = 2.9
$ws.Range("O233").Value = 3.4
$ws.Range("P233").Value = 2.4
$ws.Range("Q233").Value = 0.25
$ws.Range("R233").Value = 1.775
$ws.Range("S233").Value = 2.1
$ws.Range("T233").Value = 2.75
$ws.Range("U233").Value = 1.85
$ws.Range("V233").Value = 2
$ws.Range("W233").Value = 1.9
$ws.Range("X233").Value = -1
$ws.Range("Y233").Value = -1
$ws.Range("Z233").Value = 0.7749999999999999
$ws.Range("AA233").Value = -1
$ws.Range("AB233").Value = -1
$ws.Range("AC233").Value = 1

# Row 234
$ws.Range("B234").Value = 6078269
$ws.Range("F234").Value = "Universidad de Chile"
$ws.Range("G234").Value = "Nublense"
$ws.Range("H234").Value = 3
$ws.Range("I234").Value = 1
$ws.Range("J234").Value = "H"
$ws.Range("K234").Value = 1.85
$ws.Range("L234").Value = 3.4
$ws.Range("M234").Value = 4.333
$ws.Range("N234").Value = 1.8
$ws.Range("O234").Value = 3.6
$ws.Range("P234").Value = 4.5
$ws.Range("Q234").Value = -0.75
$ws.Range("R234").Value = 1.925
$ws.Range("S234").Value = 1.925
$ws.Range("T234").Value = 2.5
$ws.Range("U234").Value = 2.025
$ws.Range("V234").Value = 1.825
$ws.Range("W234").Value = 0.8
$ws.Range("X234").Value = -1
$ws.Range("Y234").Value = -1
$ws.Range("Z234").Value = 0.925
$ws.Range("AA234").Value = -1
$ws.Range("AB234").Value = 1.025
$ws.Range("AC234").Value = -1

# Row 235
$ws.Range("B235").Value = 6078268
$ws.Range("F235").Value = "OHiggins"
$ws.Range("G235").Value = "Palestino"
$ws.Range("H235").Value = 0
$ws.Range("I235").Value = 1
$ws.Range("J235").Value = "A"
$ws.Range("K235").Value = 3.1
$ws.Range("L235").Value = 3.3
$ws.Range("M235").Value = 2.3
$ws.Range("N235").Value = 2.9
$ws.Range("O235").Value = 3.4
$ws.Range("P235").Value = 2.375
$ws.Range("Q235").Value = 0.25
$ws.Range("R235").Value = 1.8
$ws.Range("S235").Value = 2
$ws.Range("T235").Value = 2.75
$ws.Range("U235").Value = 2
$ws.Range("V235").Value = 1.8
$ws.Range("W235").Value = -1
$ws.Range("X235").Value = -1
$ws.Range("Y235").Value = 1.375
$ws.Range("Z235").Value = -1
$ws.Range("AA235").Value = 1
$ws.Range("AB235").Value = -1
$ws.Range("AC235").Value = 0.8

# Row 236
$ws.Range("B236").Value = 6077768
$ws.Range("F236").Value = "Union La Calera"
$ws.Range("G236").Value = "Universidad Catolica"
$ws.Range("H236").Value = 0
$ws.Range("I236").Value = 3
$ws.Range("J236").Value = "A"
$ws.Range("K236").Value = 2.05
$ws.Range("L236").Value = 3.5
$ws.Range("M236").Value = 3.4
$ws.Range("N236").Value = 2.05
$ws.Range("O236").Value = 3.6
$ws.Range("P236").Value = 3.4
$ws.Range("Q236").Value = -0.25
$ws.Range("R236").Value = 1.8
$ws.Range("S236").Value = 2
$ws.Range("T236").Value = 2.75
$ws.Range("U236").Value = 1.975
$ws.Range("V236").Value = 1.825
$ws.Range("W236").Value = -1
$ws.Range("X236").Value = -1
$ws.Range("Y236").Value = 2.4
$ws.Range("Z236").Value = -1
$ws.Range("AA236").Value = 1
$ws.Range("AB236").Value = 0.4875
$ws.Range("AC236").Value = -0.5
